$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.0701390890922422
$ws.Cells.Item(2, 4).Value = 0.04032838679386685
$ws.Cells.Item(2, 5).Value = 0.07901352074717938
$ws.Cells.Item(2, 6).Value = 3.368790720019547
$ws.Cells.Item(2, 7).Value = 2.827899931590707
$ws.Cells.Item(2, 8).Value = 2.157152832125377
$ws.Cells.Item(2, 9).Value = 2.450567854647474
$ws.Cells.Item(2, 10).Value = 0.1419983326899126
$ws.Cells.Item(2, 12).Value = 0.1312270981406698
$ws.Cells.Item(2, 13).Value = 2.106278628605764
$ws.Cells.Item(2, 14).Value = 1.709010351748532

$ws.Cells.Item(3, 3).Value = 0.07019576171442665
$ws.Cells.Item(3, 4).Value = 0.0388712597791212
$ws.Cells.Item(3, 5).Value = 0.07938678922177633
$ws.Cells.Item(3, 6).Value = 3.339319399046587
$ws.Cells.Item(3, 7).Value = 2.78221808063023
$ws.Cells.Item(3, 8).Value = 2.14481270504649
$ws.Cells.Item(3, 9).Value = 2.427408465613127
$ws.Cells.Item(3, 10).Value = 0.143944590481131
$ws.Cells.Item(3, 12).Value = 0.1320584992834277
$ws.Cells.Item(3, 13).Value = 1.964616331717806
$ws.Cells.Item(3, 14).Value = 1.591420942624381

$ws.Cells.Item(4, 3).Value = 0.07023673873804981
$ws.Cells.Item(4, 4).Value = 0.03796538945656636
$ws.Cells.Item(4, 5).Value = 0.07963195359584585
$ws.Cells.Item(4, 6).Value = 3.323443527612227
$ws.Cells.Item(4, 7).Value = 2.756208700271117
$ws.Cells.Item(4, 8).Value = 2.138551629271404
$ws.Cells.Item(4, 9).Value = 2.414820296893438
$ws.Cells.Item(4, 10).Value = 0.1452005163988495
$ws.Cells.Item(4, 12).Value = 0.1326024201178626
$ws.Cells.Item(4, 13).Value = 1.877837593800123
$ws.Cells.Item(4, 14).Value = 1.519549085024892

$ws.Cells.Item(5, 3).Value = 0.07025499253766654
$ws.Cells.Item(5, 4).Value = 0.03759338453865269
$ws.Cells.Item(5, 5).Value = 0.079735888510859
$ws.Cells.Item(5, 6).Value = 3.317529384121428
$ws.Cells.Item(5, 7).Value = 2.746119209988194
$ws.Cells.Item(5, 8).Value = 2.136329729054523
$ws.Cells.Item(5, 9).Value = 2.41009875729894
$ws.Cells.Item(5, 10).Value = 0.1457276237483076
$ws.Cells.Item(5, 12).Value = 0.1328325005436479
$ws.Cells.Item(5, 13).Value = 1.842527533261972
$ws.Cells.Item(5, 14).Value = 1.4903466476338

$ws.Cells.Item(6, 3).Value = 0.07025811754684064
$ws.Cells.Item(6, 4).Value = 0.03753143972681983
$ws.Cells.Item(6, 5).Value = 0.07975339045896312
$ws.Cells.Item(6, 6).Value = 3.316580810646045
$ws.Cells.Item(6, 7).Value = 2.744474537472058
$ws.Cells.Item(6, 8).Value = 2.135980649020354
$ws.Cells.Item(6, 9).Value = 2.409339345628837
$ws.Cells.Item(6, 10).Value = 0.1458160737469942
$ws.Cells.Item(6, 12).Value = 0.1328712148635329
$ws.Cells.Item(6, 13).Value = 1.836667585926605
$ws.Cells.Item(6, 14).Value = 1.485502893696406

$ws.Cells.Item(7, 3).Value = 0.07023697861601264
$ws.Cells.Item(7, 4).Value = 0.03796038408495406
$ws.Cells.Item(7, 5).Value = 0.07963333897243885
$ws.Cells.Item(7, 6).Value = 3.323361522178061
$ws.Cells.Item(7, 7).Value = 2.756070571064441
$ws.Cells.Item(7, 8).Value = 2.138520331391504
$ws.Cells.Item(7, 9).Value = 2.414754970327166
$ws.Cells.Item(7, 10).Value = 0.145207563186585
$ws.Cells.Item(7, 12).Value = 0.1326054889102135
$ws.Cells.Item(7, 13).Value = 1.877361172699182
$ws.Cells.Item(7, 14).Value = 1.51915489837711

$ws.Cells.Item(8, 3).Value = 0.07015734797886353
$ws.Cells.Item(8, 4).Value = 0.0398282658379685
$ws.Cells.Item(8, 5).Value = 0.07913891460288003
$ws.Cells.Item(8, 6).Value = 3.358166641613494
$ws.Cells.Item(8, 7).Value = 2.811723521702191
$ws.Cells.Item(8, 8).Value = 2.152624041320564
$ws.Cells.Item(8, 9).Value = 2.442242541672812
$ws.Cells.Item(8, 10).Value = 0.1426567532166332
$ws.Cells.Item(8, 12).Value = 0.1315068375485939
$ws.Cells.Item(8, 13).Value = 2.05739262430366
$ws.Cells.Item(8, 14).Value = 1.668399276084045

$ws.Cells.Item(9, 3).Value = 0.07005017773459343
$ws.Cells.Item(9, 4).Value = 0.04340444336660454
$ws.Cells.Item(9, 5).Value = 0.07829561255872886
$ws.Cells.Item(9, 6).Value = 3.444164537521033
$ws.Cells.Item(9, 7).Value = 2.937200611373157
$ws.Cells.Item(9, 8).Value = 2.19078702994014
$ws.Cells.Item(9, 9).Value = 2.509194395350875
$ws.Cells.Item(9, 10).Value = 0.1381383500697555
$ws.Cells.Item(9, 12).Value = 0.1296167824860959
$ws.Cells.Item(9, 13).Value = 2.411969417942259
$ws.Cells.Item(9, 14).Value = 1.963545378273579

$ws.Cells.Item(10, 3).Value = 0.07000125405009783
$ws.Cells.Item(10, 4).Value = 0.04598202965242848
$ws.Cells.Item(10, 5).Value = 0.07775232391770004
$ws.Cells.Item(10, 6).Value = 3.518363652810734
$ws.Cells.Item(10, 7).Value = 3.039587960743233
$ws.Cells.Item(10, 8).Value = 2.225327803561498
$ws.Cells.Item(10, 9).Value = 2.566489567977229
$ws.Cells.Item(10, 10).Value = 0.1351139928571525
$ws.Cells.Item(10, 12).Value = 0.1283880808755029
$ws.Cells.Item(10, 13).Value = 2.673346933605814
$ws.Cells.Item(10, 14).Value = 2.181753477840402

$ws.Cells.Item(11, 3).Value = 0.06998546366672898
$ws.Cells.Item(11, 4).Value = 0.04714449631204332
$ws.Cells.Item(11, 5).Value = 0.07752158597416781
$ws.Cells.Item(11, 6).Value = 3.5545553592587
$ws.Cells.Item(11, 7).Value = 3.088434233830696
$ws.Cells.Item(11, 8).Value = 2.242475271690466
$ws.Cells.Item(11, 9).Value = 2.594348453704825
$ws.Cells.Item(11, 10).Value = 0.1338023513148801
$ws.Cells.Item(11, 12).Value = 0.1278635702591586
$ws.Cells.Item(11, 13).Value = 2.792431668994141
$ws.Cells.Item(11, 14).Value = 2.281289632159201

$ws.Cells.Item(12, 3).Value = 0.06998041327729965
$ws.Cells.Item(12, 4).Value = 0.04758329241521153
$ws.Cells.Item(12, 5).Value = 0.07743655926504012
$ws.Cells.Item(12, 6).Value = 3.568614296144943
$ws.Cells.Item(12, 7).Value = 3.107261556160438
$ws.Cells.Item(12, 8).Value = 2.249176600339069
$ws.Cells.Item(12, 9).Value = 2.605158648384148
$ws.Cells.Item(12, 10).Value = 0.1333149044602218
$ws.Cells.Item(12, 12).Value = 0.1276698823614204
$ws.Cells.Item(12, 13).Value = 2.837550734306603
$ws.Cells.Item(12, 14).Value = 2.319017725263848

$ws.Cells.Item(13, 3).Value = 0.06998145965727787
$ws.Cells.Item(13, 4).Value = 0.0474888515591374
$ws.Cells.Item(13, 5).Value = 0.07745476699468323
$ws.Cells.Item(13, 6).Value = 3.56557065817276
$ws.Cells.Item(13, 7).Value = 3.103192002280366
$ws.Cells.Item(13, 8).Value = 2.247724073867175
$ws.Cells.Item(13, 9).Value = 2.60281884694281
$ws.Cells.Item(13, 10).Value = 0.1334194734687193
$ws.Cells.Item(13, 12).Value = 0.12771137741559
$ws.Cells.Item(13, 13).Value = 2.827832489045534
$ws.Cells.Item(13, 14).Value = 2.310890758972334

$ws.Cells.Item(14, 3).Value = 0.06998502955219266
$ws.Cells.Item(14, 4).Value = 0.0471806242189885
$ws.Cells.Item(14, 5).Value = 0.07751454376359579
$ws.Cells.Item(14, 6).Value = 3.555704884077102
$ws.Cells.Item(14, 7).Value = 3.089976526642175
$ws.Cells.Item(14, 8).Value = 2.243022416757526
$ws.Cells.Item(14, 9).Value = 2.595232576743612
$ws.Cells.Item(14, 10).Value = 0.1337620634136019
$ws.Cells.Item(14, 12).Value = 0.1278475366812799
$ws.Cells.Item(14, 13).Value = 2.796143166865505
$ws.Cells.Item(14, 14).Value = 2.284392843647481

$ws.Cells.Item(15, 3).Value = 0.06998733719325756
$ws.Cells.Item(15, 4).Value = 0.04699164451082538
$ws.Cells.Item(15, 5).Value = 0.07755146432369919
$ws.Cells.Item(15, 6).Value = 3.549708003391515
$ws.Cells.Item(15, 7).Value = 3.081924798132434
$ws.Cells.Item(15, 8).Value = 2.240169647310552
$ws.Cells.Item(15, 9).Value = 2.590619785630665
$ws.Cells.Item(15, 10).Value = 0.1339731139455096
$ws.Cells.Item(15, 12).Value = 0.1279315800272016
$ws.Cells.Item(15, 13).Value = 2.77673564216667
$ws.Cells.Item(15, 14).Value = 2.268166687334997

$ws.Cells.Item(16, 3).Value = 0.07000241603268975
$ws.Cells.Item(16, 4).Value = 0.04590586061239321
$ws.Cells.Item(16, 5).Value = 0.07776773243066337
$ws.Cells.Item(16, 6).Value = 3.516047783921579
$ws.Cells.Item(16, 7).Value = 3.036441747922652
$ws.Cells.Item(16, 8).Value = 2.224236185962468
$ws.Cells.Item(16, 9).Value = 2.564705258833087
$ws.Cells.Item(16, 10).Value = 0.1352010041917722
$ws.Cells.Item(16, 12).Value = 0.128423050163164
$ws.Cells.Item(16, 13).Value = 2.665567970168041
$ws.Cells.Item(16, 14).Value = 2.1752537558732

$ws.Cells.Item(17, 3).Value = 0.07001332185451048
$ws.Cells.Item(17, 4).Value = 0.04523721562673444
$ws.Cells.Item(17, 5).Value = 0.07790460058114768
$ws.Cells.Item(17, 6).Value = 3.496025141741057
$ws.Cells.Item(17, 7).Value = 3.009123469261851
$ws.Cells.Item(17, 8).Value = 2.214830144773089
$ws.Cells.Item(17, 9).Value = 2.54926906073571
$ws.Cells.Item(17, 10).Value = 0.135970719139884
$ws.Cells.Item(17, 12).Value = 0.1287333569056557
$ws.Cells.Item(17, 13).Value = 2.597415565979475
$ws.Cells.Item(17, 14).Value = 2.118322080474172

$ws.Cells.Item(18, 3).Value = 0.07002020316499191
$ws.Cells.Item(18, 4).Value = 0.04485167616834218
$ws.Cells.Item(18, 5).Value = 0.07798486843478347
$ws.Cells.Item(18, 6).Value = 3.4847379361712
$ws.Cells.Item(18, 7).Value = 2.993624119955086
$ws.Cells.Item(18, 8).Value = 2.20955497860416
$ws.Cells.Item(18, 9).Value = 2.540559370723685
$ws.Cells.Item(18, 10).Value = 0.1364194761661515
$ws.Cells.Item(18, 12).Value = 0.1289150791254272
$ws.Cells.Item(18, 13).Value = 2.558233429035653
$ws.Cells.Item(18, 14).Value = 2.085602296101456

$ws.Cells.Item(19, 3).Value = 0.07002263760107397
$ws.Cells.Item(19, 4).Value = 0.04472097413574261
$ws.Cells.Item(19, 5).Value = 0.07801231141212828
$ws.Cells.Item(19, 6).Value = 3.480955562702633
$ws.Cells.Item(19, 7).Value = 2.988412840534352
$ws.Cells.Item(19, 8).Value = 2.207792024392973
$ws.Cells.Item(19, 9).Value = 2.537639335444851
$ws.Cells.Item(19, 10).Value = 0.1365724539580842
$ws.Cells.Item(19, 12).Value = 0.1289771645217499
$ws.Cells.Item(19, 13).Value = 2.54497006848618
$ws.Cells.Item(19, 14).Value = 2.074528471526037

$ws.Cells.Item(20, 3).Value = 0.07001209793045149
$ws.Cells.Item(20, 4).Value = 0.04530849227999312
$ws.Cells.Item(20, 5).Value = 0.07788987091661159
$ws.Cells.Item(20, 6).Value = 3.498132831843918
$ws.Cells.Item(20, 7).Value = 3.012009432256548
$ws.Cells.Item(20, 8).Value = 2.215817456328239
$ws.Cells.Item(20, 9).Value = 2.550894781843255
$ws.Cells.Item(20, 10).Value = 0.1358881566978525
$ws.Cells.Item(20, 12).Value = 0.1286999888233993
$ws.Cells.Item(20, 13).Value = 2.604668727662926
$ws.Cells.Item(20, 14).Value = 2.12437990690745

$ws.Cells.Item(21, 3).Value = 0.06998395577587857
$ws.Cells.Item(21, 4).Value = 0.04727119581429662
$ws.Cells.Item(21, 5).Value = 0.07749692221725546
$ws.Cells.Item(21, 6).Value = 3.558593068938279
$ws.Cells.Item(21, 7).Value = 3.093849236455696
$ws.Cells.Item(21, 8).Value = 2.244397750292194
$ws.Cells.Item(21, 9).Value = 2.597453754222869
$ws.Cells.Item(21, 10).Value = 0.1336611854342844
$ws.Cells.Item(21, 12).Value = 0.1278074096733555
$ws.Cells.Item(21, 13).Value = 2.805450446427926
$ws.Cells.Item(21, 14).Value = 2.292174974504576

$ws.Cells.Item(22, 3).Value = 0.06997097846754841
$ws.Cells.Item(22, 4).Value = 0.04854576978338798
$ws.Cells.Item(22, 5).Value = 0.07725379291625867
$ws.Cells.Item(22, 6).Value = 3.600171849385703
$ws.Cells.Item(22, 7).Value = 3.149263378523244
$ws.Cells.Item(22, 8).Value = 2.264289532260079
$ws.Cells.Item(22, 9).Value = 2.629403141433869
$ws.Cells.Item(22, 10).Value = 0.1322596168824073
$ws.Cells.Item(22, 12).Value = 0.1272528017879324
$ws.Cells.Item(22, 13).Value = 2.936813409456079
$ws.Cells.Item(22, 14).Value = 2.40204661557118

$ws.Cells.Item(23, 3).Value = 0.06997740938394159
$ws.Cells.Item(23, 4).Value = 0.04786623743321883
$ws.Cells.Item(23, 5).Value = 0.07738230686221037
$ws.Cells.Item(23, 6).Value = 3.57779045475138
$ws.Cells.Item(23, 7).Value = 3.119510167939154
$ws.Cells.Item(23, 8).Value = 2.25356137950007
$ws.Cells.Item(23, 9).Value = 2.612211190739231
$ws.Cells.Item(23, 10).Value = 0.1330027237399087
$ws.Cells.Item(23, 12).Value = 0.1275461824080537
$ws.Cells.Item(23, 13).Value = 2.866690324455789
$ws.Cells.Item(23, 14).Value = 2.343388084879905

$ws.Cells.Item(24, 3).Value = 0.07001264936251061
$ws.Cells.Item(24, 4).Value = 0.04527627161787251
$ws.Cells.Item(24, 5).Value = 0.07789652527402247
$ws.Cells.Item(24, 6).Value = 3.497179247445047
$ws.Cells.Item(24, 7).Value = 3.01070404612463
$ws.Cells.Item(24, 8).Value = 2.215370680214136
$ws.Cells.Item(24, 9).Value = 2.550159280248877
$ws.Cells.Item(24, 10).Value = 0.1359254637500786
$ws.Cells.Item(24, 12).Value = 0.1287150641825079
$ws.Cells.Item(24, 13).Value = 2.601389574394148
$ws.Cells.Item(24, 14).Value = 2.121641129086584

$ws.Cells.Item(25, 3).Value = 0.07007393122528072
$ws.Cells.Item(25, 4).Value = 0.04244596567370706
$ws.Cells.Item(25, 5).Value = 0.07851030196583686
$ws.Cells.Item(25, 6).Value = 3.418979775877375
$ws.Cells.Item(25, 7).Value = 2.901481971012288
$ws.Cells.Item(25, 8).Value = 2.179328224270392
$ws.Cells.Item(25, 9).Value = 2.489669757359565
$ws.Cells.Item(25, 10).Value = 0.1393088695637883
$ws.Cells.Item(25, 12).Value = 0.1300999177186348
$ws.Cells.Item(25, 13).Value = 2.315889609796187
$ws.Cells.Item(25, 14).Value = 1.883450603060282

